$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. "496.49") are not silently converted to numbers, matching
# the source data which stores every Price/Volume cell as text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '56.784.04'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '2.964.31'
$ws.Range("E3").Value = '  -1.56%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '496.49'
$ws.Range("E5").Value = '  -3.64%  '
$ws.Range("D6").Value = '137.29'
$ws.Range("E6").Value = '  -1.91%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '0.425'
$ws.Range("E8").Value = '  -2.50%  '
$ws.Range("D9").Value = '7.30'
$ws.Range("E10").Value = '  -2.52%  '
$ws.Range("D11").Value = '0.356'
$ws.Range("E11").Value = '  -0.75%  '
$ws.Range("D12").Value = '3.469.65'
$ws.Range("E12").Value = '  -1.62%  '
$ws.Range("D13").Value = '0.128'
$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("D14").Value = '25.73'
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("D16").Value = '56.827.52'
$ws.Range("E16").Value = '  -0.42%  '
$ws.Range("D17").Value = '6.05'
$ws.Range("E17").Value = '  +1.66%  '
$ws.Range("D18").Value = '2.959.17'
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("D19").Value = '12.55'
$ws.Range("E19").Value = '  -0.59%  '
$ws.Range("D20").Value = '7.78'
$ws.Range("E20").Value = '  -1.70%  '
$ws.Range("D21").Value = '318.28'
$ws.Range("E21").Value = '  -3.33%  '
$ws.Range("E22").Value = '  -0.38%  '
$ws.Range("D23").Value = '5.64'
$ws.Range("E23").Value = '  -0.90%  '
$ws.Range("E24").Value = '  -0.35%  '
$ws.Range("D25").Value = '63.09'
$ws.Range("E25").Value = '  -0.98%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = '0.162'
$ws.Range("E27").Value = '  -5.75%  '
$ws.Range("D28").Value = '0.0₃0885'
$ws.Range("E28").Value = '  -3.75%  '
$ws.Range("D29").Value = '6.51'
$ws.Range("E29").Value = '  -2.28%  '
$ws.Range("D30").Value = '7.03'
$ws.Range("E30").Value = '  -1.67%  '
$ws.Range("E31").Value = '  -3.59%  '
$ws.Range("E32").Value = '  -6.34%  '
$ws.Range("D33").Value = '20.05'
$ws.Range("D34").Value = '155.54'
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("D35").Value = '4.60'
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '5.71'
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("D37").Value = '1.25'
$ws.Range("E37").Value = '  -3.03%  '
$ws.Range("D38").Value = '23.85'
$ws.Range("D39").Value = '0.0663'
$ws.Range("E39").Value = '  -2.63%  '
$ws.Range("D40").Value = '2.994.67'
$ws.Range("E40").Value = '  -1.62%  '
$ws.Range("D41").Value = '37.31'
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("D43").Value = '3.70'
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("E44").Value = '  -2.15%  '
$ws.Range("D45").Value = '2.193.98'
$ws.Range("E45").Value = '  -4.44%  '
$ws.Range("E46").Value = '  -3.80%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").Value = '0.934'
$ws.Range("E47").Value = '  -7.32%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '5.90'
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("D49").Value = '0.0234'
$ws.Range("E49").Value = '  -3.16%  '
$ws.Range("D50").Value = '19.14'
$ws.Range("E50").Value = '  -1.09%  '
$ws.Range("D51").Value = '1.78'
$ws.Range("E51").Value = '  -10.50%  '

# Restore the original (default) cell style on column D now that the
# text values are committed, so no residual number-format styling remains.
$ws.Range("D2:D51").Style = "Normal"
